$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlTop = -4160
$xlPasteFormats = -4122

# ---- Row 1 header: add column G "Missing" (copy header look from C1) ----
$ws.Range("C1").Copy()
$ws.Range("G1").PasteSpecial($xlPasteFormats)
$ws.Range("G1").Value = 'Missing'

# ---- B2, B3: new criteria text; copy bold/wrap/top format from B1 ----
$ws.Range("B1").Copy()
$ws.Range("B2").PasteSpecial($xlPasteFormats)
$ws.Range("B3").PasteSpecial($xlPasteFormats)
$ws.Range("B2").Value = 'Clarity of Big Idea / Primary Use Case for a Non-Technical Audience'
$ws.Range("B3").Value = 'Tension Development, Detail and Persuasiveness of  Argumentation in Data Story'

# ---- B4: new "Effectiveness of Pitch" label; keep AptosDisplay bold font, add vertical-top ----
$ws.Range("B4").Value = 'Effectiveness of Pitch'
$ws.Range("B4").VerticalAlignment = $xlTop

# ---- B5: new criteria text only (style already correct: bold/top/wrap) ----
$ws.Range("B5").Value = 'Quality and Relevance of Visualizations'

# ---- C2:F5 + G2:G6: new descriptive text; copy default font + vertical-top + wrap look from A2 ----
$ws.Range("A2").Copy()
$ws.Range("C2:F5").PasteSpecial($xlPasteFormats)
$ws.Range("G2:G6").PasteSpecial($xlPasteFormats)
$ws.Range("C2").Value = 'Visually through smart use of text placement, spacing, style, visuals and audibly through the team''s speakers'' tone, pace, and other oratory techniques the audience is guided effectively and convincingly to the final selected use case, through the data story and to the final pitch.'
$ws.Range("D2").Value = 'Verbal and visual presentation are aligned well to communicate the Big Idea / primary use case effectively. The final selected use case is motivasted convincingly and the data story is told in an overall engaging manner to end in a well-argued pitch.'
$ws.Range("E2").Value = 'The final selected use case is mentioned but could have been better separated from other use cases mentioned. The story-telling with the data feels disjointed or does not sufficiently engage the audience to create a convincing final pitch.'
$ws.Range("F2").Value = 'It is not clear from the presentation what the final selected use case is. The direction and ultimate objectives of the data story are not clear. The final pitch to the stakeholders is not clear.'
$ws.Range("C3").Value = 'Visual and verbal presentation align perfectly to elegantly and convincingly walk the audience through the business problem, the use case selection and the data story. The reasoning is persuasive and well designed to also convince non-technical business stakeholders. The data exploration illustrates a high level of thought and insight with a feasibility study for use cases clearly outlined and motivated. The logic succession of steps in the data processing is easily understandable.'
$ws.Range("D3").Value = 'Visual and verbal presentation components integrate well to make the case for the selection of steps in the data exploration process. The motivation is easy to follow and is convincing based on the data and the analysis presented. The data exploration is clearly described and well-formulated as a data story. The logic in the data exploration can be understood and conclusions make sense and carry the data story forward to motivate the proposed solution for the business problem.'
$ws.Range("E3").Value = 'The presentation shows all the relevant steps in the data exploration process. But the motivation for the selection and prioritization of steps In the data exploration process is not clear. A data exploration process is described with some explanation of how certain aspects in the data will affect the use case discovery process. But the description of feature importance and feature correlations and correlation to the target feature is not clear and it is hard to follow the data story towards a solution of the business problem.'
$ws.Range("F3").Value = 'It is not clear how some or any of the steps in the data exploration process are relevant for the data science use case and the business problem. The data exploration steps are not clearly outlined or missing. The process to explore different use cases for the business problem at hand is not clearly explained.'
$ws.Range("C4").Value = 'A well-crafted story emerges from the presentation of how the proposed solution will address the customer''s business problem in a language that a non-technical stakeholder will understand. Tone, pace and delivery are well adjusted to communicate the pitch effectively in the total 5-7 min time limit allotted to the whole presentation. In the presentation it can be clearly seen that the pitch is conducted with the customer''s perspective in mind. The solution presented is clearly solving their problem, or an aspect of it, and the audience can convincingly see the business value in the solution. A call to action very nicely emerges from the flow of the presentation.'
$ws.Range("D4").Value = 'An explanation is provided of how the service/product/solution you are proposing will solve for the problem. Points flow logically. Usage scenarios are described and how you will address them in simple, customer-understood language.  A customer journey is painted from the perspective of the customer where it is clear from the description how your team will discover and apply the solution. You are able to outline your value proposition, i.e. how the customer will react / feel about your solution. A call to action is formulated and can be directly derived from the explanations on the above items. There are no questions remaining of what is requested from the stakeholders.'
$ws.Range("E4").Value = 'The business problem and your team''s solution are connected and this can be seen from your presentation. But the argumentation does not flow easily, key questions remain whether a better approach exists or whether your solution would truly create business value for the stakeholder. The business problem and your team''s solution are connected and this can be seen from your presentation. A call to action is formulated but could be made more precise and some questions remain.'
$ws.Range("F4").Value = 'There are logical gaps in how your solution connects with the business problem of your customer and the description in the presentation is not easy to follow. There is a disconnect between the customer''s problem(s) and the solution your team is proposing. There is the impression that the presented problem serves to justify your team''s solution and not the opposite. The business value of your solution is not convincingly presented to the audience. A call to action for the stakeholders is not clearly formulated and confusing.'
$ws.Range("C5").Value = 'The visualizations stand out in terms of quality, style, clarity. They are able to help focus the key messages in the data story. Color, graph choice, labeling, descriptions are thoughtfully and effectively used to create high-quality professional visualizations. With the maximally 4 visualizations the audience is taken on a journey that leads through the essential aspects of the problem, the use case and the data to persuasively pitch the solution proposed.'
$ws.Range("D5").Value = 'The visualizations used were relevant and informative and logically created a path for the audience to understand the particular complexities of the data and how handling those in logical sequence will lead to an effective solution of the business problem. This is achieved with the help of maximally 4 well-crafted visuals tied into a description of business problem and use case and the data story and concludes with a data science solution proposal of the business problem. Space, color, the choice of graphing style and other elements add effectively to the telling of the overall story,'
$ws.Range("E5").Value = 'Three visualizations were added to the data story and one to the business problem / use case description. They were useful in communicating essential parts of the data story but visuals or their descriptions in the presentation could have been improved to tell the data story more clearly and connect more effectively to the final pitch to stakeholders.'
$ws.Range("F5").Value = 'Visualizations were included in the slide deck but only with limited relevance to telling the data story, business problem or data science use case of how the data exploration can help solve the business problem presented. Key elements in the graph like axis labels are missing or hard to read or generally visualizations are missing and/or of lower quality.'
$ws.Range("G2").Value = 'Missing or no work was submitted.'
$ws.Range("G3").Value = 'Missing or no work was submitted.'
$ws.Range("G4").Value = 'Missing or no work was submitted.'
$ws.Range("G5").Value = 'Missing or no work was submitted.'
$ws.Range("G6").Value = 'Missing or no work was submitted.'

# ---- C6:F6: new text, existing filler style already correct (no format change) ----
$ws.Range("C6").Value = 'The PowerPoint presentation is exceptional in design and delivery. Slides are cohesive with a consistent design, using appropriate fonts, colors, and objects. The presentation flows smoothly, is engaging, and keeps the audience''s attention. The presenter delivers the content confidently and clearly within the allotted time slot. '
$ws.Range("D6").Value = 'The PowerPoint presentation is well- designed and delivered. Slides exhibit a good degree of cohesion in terms of design elements. The presentation is clear and mostly fluent, and the presenter manages to stay within the allotted time. '
$ws.Range("E6").Value = 'The PowerPoint presentation lacks some cohesiveness in design elements, making it less visually appealing. The delivery is somewhat choppy but generally understandable, and the presentation mostly fits the allotted time. '
$ws.Range("F6").Value = 'The PowerPoint presentation is poorly designed, lacking consistency in fonts, colors, and objects. The delivery is disjointed or difficult to follow, and the presentation significantly exceeds or falls short of the allotted time. '

# ---- Row heights to match content ----
$ws.Rows.Item(2).RowHeight = 102
$ws.Rows.Item(3).RowHeight = 187
$ws.Rows.Item(4).RowHeight = 255
$ws.Rows.Item(5).RowHeight = 221
$ws.Rows.Item(6).RowHeight = 119

# ---- Selection ----
$ws.Range("D4").Select()
